# Add a new login row ("Admin3" / "admin123") to the LoginData sheet,
# mirroring the existing Admin/Admin1 rows (row 3, columns A:B).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

$ws.Range("A3").Value = "Admin3"
$ws.Range("B3").Value = "admin123"

# Restore the view to the top-left cell (the prior saved selection pointed
# at B6, which no longer reflects useful UI state after this edit).
[void]$ws.Range("A1").Select()
